$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab year from 2025 to 2023
$ws.Name = "g3.5a 2023"

# Update region order (rows 3-7 reordered) and refreshed values (column B)
$ws.Range("B2").Value = 32.069953084778

$ws.Range("A3").Value = "Mato Grosso do Sul"
$ws.Range("B3").Value = 17.53196975891699

$ws.Range("A4").Value = "Mato Grosso"
$ws.Range("B4").Value = 17.49915955766788

$ws.Range("A5").Value = "Sergipe"
$ws.Range("B5").Value = 17.46029783768665

$ws.Range("B6").Value = 15.89086868026616

$ws.Range("A7").Value = "Goiás"
$ws.Range("B7").Value = 14.91689084958633

$ws.Range("B8").Value = 14.29087477237749

$ws.Range("B9").Value = 8.045879236039708
